# Added "optimal parameters algorithm" block: two new mini-tables
# (rows 18-22 and rows 24-28) mirroring the existing efficiency tables,
# plus an updated selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell, far outside the used range, used purely to stamp a
# "Text" number format onto the percentage cells before their values are
# assigned - this keeps "0.00008%"-style strings from being reinterpreted
# as numeric percentages. It is cleared again at the end.
$stamp = $ws.Range("Z100")
$stamp.NumberFormat = "@"

# ---------------------------------------------------------------------
# Step 1: write all the new cell values first, in the same order the
# strings end up in the shared-string table.
# ---------------------------------------------------------------------

$stamp.Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("I20").Value = "0.00002%"

$ws.Range("H18").Value = "Muon1.pt > 17.5,  Photon.pt > 23.5, 2.25<m(JPsi)<3.5"

$stamp.Copy()
$ws.Range("I21").PasteSpecial(-4122)
$ws.Range("I21").Value = "25.13%"

$ws.Range("H24").Value = "Muon1.pt > 15.,  Photon.pt > 23.5, 2.7<m(JPsi)<3.5"

$stamp.Copy()
$ws.Range("I26").PasteSpecial(-4122)
$ws.Range("I26").Value = "0.00005%"

$stamp.Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("I27").Value = "26.31%"

$stamp.Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I28").Value = "37.68%"

$stamp.Copy()
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("I22").Value = "37.61%"

$ws.Range("F18").Value = "Muon1.pt > 20, Photon.pt > 26.3, dR(mumu) < 0.35"

$stamp.Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value = "0.00008%"

$stamp.Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("G21").Value = "22.34%"

$stamp.Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = "36.62%"

$ws.Range("G19").Value = "selection"
$ws.Range("I19").Value = "selection"
$ws.Range("I25").Value = "selection"

$stamp.Clear()

# ---------------------------------------------------------------------
# Step 2: merge the header ranges and apply formatting that mirrors the
# existing tables (thin box border + bold for headers, plain centered
# style for the data rows). The values already written in Step 1 are
# left untouched so they stay stored as literal text.
# ---------------------------------------------------------------------

$ws.Range("F18:G18").Merge()
$ws.Range("F12:G12").Copy()
$ws.Range("F18:G18").PasteSpecial(-4122)

$ws.Range("H18:I18").Merge()
$ws.Range("H12:I12").Copy()
$ws.Range("H18:I18").PasteSpecial(-4122)

$ws.Range("H24:I24").Merge()
$ws.Range("H12:I12").Copy()
$ws.Range("H24:I24").PasteSpecial(-4122)

$ws.Range("G13").Copy()
$ws.Range("G19").PasteSpecial(-4122)

$ws.Range("I13").Copy()
$ws.Range("I19").PasteSpecial(-4122)

$ws.Range("I13").Copy()
$ws.Range("I25").PasteSpecial(-4122)

$ws.Range("G14").Copy()
$ws.Range("G20").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("I20").PasteSpecial(-4122)

$ws.Range("G15").Copy()
$ws.Range("G21").PasteSpecial(-4122)

$ws.Range("I15").Copy()
$ws.Range("I21").PasteSpecial(-4122)

$ws.Range("G16").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("I16").Copy()
$ws.Range("I22").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("I26").PasteSpecial(-4122)

$ws.Range("I15").Copy()
$ws.Range("I27").PasteSpecial(-4122)

$ws.Range("I16").Copy()
$ws.Range("I28").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 3: view state - active selection
# ---------------------------------------------------------------------
$ws.Range("G23").Select() | Out-Null
